# Weekly update: a new price record (Packham's Triumph, Calibre 70,
# Provincia de Curico) is inserted as a new row 543 in the "Pera" sheet.
# All the existing records previously occupying rows 543-574 shift down
# by one row, so the sheet grows from A1:T574 to A1:T575.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 543; Excel shifts rows 543:574 -> 544:575
# and keeps per-column formatting (e.g. the date style on column D).
$ws.Rows.Item(543).Insert()

$curico = "Provincia de Curic" + [char]0x00F3

# Populate the newly-inserted row 543 with the new record.
$ws.Cells.Item(543, 1).Value  = 9
$ws.Cells.Item(543, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(543, 3).Value  = "Metropolitana"
$ws.Cells.Item(543, 4).Value  = 44509
$ws.Cells.Item(543, 5).Value  = 13
$ws.Cells.Item(543, 6).Value  = "Fruta"
$ws.Cells.Item(543, 7).Value  = 100104
$ws.Cells.Item(543, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(543, 9).Value  = 100104005
$ws.Cells.Item(543, 10).Value = "Pera"
$ws.Cells.Item(543, 11).Value = "Packham's Triumph"
$ws.Cells.Item(543, 12).Value = "Calibre 70"
$ws.Cells.Item(543, 13).Value = 570
$ws.Cells.Item(543, 14).Value = 18000
$ws.Cells.Item(543, 15).Value = 19000
$ws.Cells.Item(543, 16).Value = 18614
$ws.Cells.Item(543, 17).Value = "`$/caja 18 kilos embalada"
$ws.Cells.Item(543, 18).Value = $curico
$ws.Cells.Item(543, 19).Value = 1034
$ws.Cells.Item(543, 20).Value = 18
